# Update "想去人数" (want-to-go count) values in both the "展览" and
# "全部类型" worksheets, which hold duplicate copies of the same data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 1371
    "F3" = 2115
    "F4" = 292
    "F6" = 6383
    "F7" = 265
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
